$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tweak wording on the two most-recent existing tasting notes (row 14 "pop", row 15 "small puff") ---
# Edit D15 first, then D14, to reproduce the shared-string slot order seen in the target file.
$ws.Range("D15").Value2 = "Cooled and served @ 12 C. Opening the swing-top gave a small puff. Poured clear with foam. Moderate carbonation. A light sour grassy taste, less boozy and a dry mouthfeel."
$ws.Range("D14").Value2 = "Cooled and served @ 12 C. Opening the swing-top gave a pop. Poured clear with foam. Moderate carbonation. A light sour grassy taste, less boozy and a dry mouthfeel."

# --- Add the new tasting-note row (row 16) ---
$ws.Range("A16").Value2 = 44204
$ws.Range("B16").Formula = '=A16-$A$6'
$ws.Range("C16").Value2 = 2.75
$ws.Range("D16").Value2 = "Cooled and served @ 12 C. Opening the swing-top gave a small puff. Poured clear with a little foam. Moderate carbonation. A light sour grassy taste, less boozy and a dry mouthfeel."

# Copy the formatting of the row above down onto the new row (border, number format, alignment, wrap).
$ws.Range("A15:D15").Copy() | Out-Null
$ws.Range("A16:D16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(15).RowHeight

# --- View state: move selection to the new next-empty row ---
$ws.Range("D17").Select() | Out-Null

# --- Window geometry (best effort) ---
$excel.ActiveWindow.Left = 20370
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 25440
$excel.ActiveWindow.Height = 15540

Write-Host "done"
